$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "imageUrl:" -> "eventImage" + ":" (two runs, same Menlo/9CDCFE formatting)
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("imageUrl:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$fullStart = $r.Start
$fullEnd = $r.End

# Replace "imageUrl" (everything but the trailing colon) with "eventImage"
$labelPart = $d.Range($fullStart, $fullEnd - 1)
$labelPart.Text = "eventImage"

# The trailing colon keeps its own run; force Word to split it into a
# distinct <w:r> (identical formatting to the label) by toggling a
# character property on/off instead of leaving it merged with the label run.
$newEnd = $labelPart.End
$colonPart = $d.Range($newEnd, $newEnd + 1)
$colonPart.Bold = $true
$colonPart2 = $d.Range($newEnd, $newEnd + 1)
$colonPart2.Bold = $false

# ---------------------------------------------------------------------------
# 2) Append two new paragraphs to the ADD/EDIT EVENT form at the end of the
#    document: one with two spaces, one with an em dash in Calibri.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$spacesPara = $d.Paragraphs.Last
$spacesPara.Range.Text = "  "

$dashPara = $d.Paragraphs.Last
$dashPara.Range.InsertParagraphAfter()

$emDashPara = $d.Paragraphs.Last
$emDashPara.Range.Text = "—"
$dashRange = $d.Range($emDashPara.Range.Start, $emDashPara.Range.Start + 1)
$dashRange.Font.Name = "Calibri"

Write-Output "done"
